# Updates cryptos list values (Price column D, Volume(1h) column E)
# per commit "Updated cryptos list on Sun Oct 15 07:24:06 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.008.35"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.46"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.28"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.06"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.783.60"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.560.60"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.023.61"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.81"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.49"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.21"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.421.80"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("E36").Value = "  +10.22%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.65"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.697.60"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.61"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  +0.67%  "
